# This edit swaps the data block occupying rows 2-3 with the data block
# occupying rows 4-5 on the active worksheet (a cyclic re-ordering of the
# four data rows: new2=old4, new3=old5, new4=old2, new5=old3).
#
# A temporary helper worksheet is used as scratch space so the swap can be
# done with simple Range.Copy operations (which - unlike assigning through
# .Value/.Value2 - preserve the original cell types/contents verbatim,
# instead of letting Excel "smart" re-parse text that looks like a
# date/number).

$wb = $excel.ActiveWorkbook
$mainName = $wb.ActiveSheet.Name
$ws = $wb.Worksheets.Item($mainName)

$dataCols = "AY"   # last used column in the sheet
$blockTop = "A2:" + $dataCols + "3"
$blockBottom = "A4:" + $dataCols + "5"
$tempBlock = "A1:" + $dataCols + "2"

# Create a scratch worksheet to stage one of the blocks during the swap.
$tempSheet = $wb.Worksheets.Add()
$tempName = $tempSheet.Name

# Re-fetch sheet references by name (ActiveSheet-bound references can shift
# when sheets are added/removed).
$ws = $wb.Worksheets.Item($mainName)
$tempSheet = $wb.Worksheets.Item($tempName)

# 1) Stage rows 2:3 on the temp sheet.
$tempSheet.Range($tempBlock).ClearContents()
$ws.Range($blockTop).Copy($tempSheet.Range($tempBlock))

# 2) Move rows 4:5 into rows 2:3 (clear destination first so cells that are
#    blank in the source truly become blank, rather than retaining stale
#    destination content).
$ws.Range($blockTop).ClearContents()
$ws.Range($blockBottom).Copy($ws.Range($blockTop))

# 3) Move the staged (original rows 2:3) content into rows 4:5.
$ws.Range($blockBottom).ClearContents()
$tempSheet.Range($tempBlock).Copy($ws.Range($blockBottom))

# Clean up the scratch worksheet.
$tempSheet.Delete()

# Re-fetch the main worksheet reference once more (it becomes stale as soon
# as another sheet is deleted).
$ws = $wb.Worksheets.Item($mainName)
